$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 11 (existing rows 12..26 shift down to 14..28)
$ws.Rows("12:13").Insert()

# New row 12 data
$ws.Range("A12").Value2 = 11
$ws.Range("B12").Value2 = "Vega Monumental Concepción"
$ws.Range("C12").Value2 = "Bíobío"
$ws.Range("D12").Value2 = 44483
$ws.Range("E12").Value2 = 8
$ws.Range("F12").Value2 = 100112013
$ws.Range("G12").Value2 = "Alcachofa"
$ws.Range("H12").Value2 = "Española"
$ws.Range("I12").Value2 = "Primera"
$ws.Range("J12").Value2 = 450
$ws.Range("K12").Value2 = 11000
$ws.Range("L12").Value2 = 12000
$ws.Range("M12").Value2 = 11444
$ws.Range("N12").Value2 = "$/caja 30 unidades"
$ws.Range("O12").Value2 = "Provincia de Limarí"
$ws.Range("P12").Value2 = 381
$ws.Range("Q12").Value2 = 30
$ws.Range("R12").Value2 = "Hortaliza"

# New row 13 data
$ws.Range("A13").Value2 = 11
$ws.Range("B13").Value2 = "Vega Monumental Concepción"
$ws.Range("C13").Value2 = "Bíobío"
$ws.Range("D13").Value2 = 44483
$ws.Range("E13").Value2 = 8
$ws.Range("F13").Value2 = 100112013
$ws.Range("G13").Value2 = "Alcachofa"
$ws.Range("H13").Value2 = "Madrigal"
$ws.Range("I13").Value2 = "Primera"
$ws.Range("J13").Value2 = 220
$ws.Range("K13").Value2 = 8000
$ws.Range("L13").Value2 = 8500
$ws.Range("M13").Value2 = 8273
$ws.Range("N13").Value2 = "$/caja 40 unidades"
$ws.Range("O13").Value2 = "Región de Coquimbo"
$ws.Range("P13").Value2 = 207
$ws.Range("Q13").Value2 = 40
$ws.Range("R13").Value2 = "Hortaliza"
